$wb = $excel.ActiveWorkbook

# ---- Worksheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 42
$ws.Range("H42").Value = 634.55554
$ws.Range("I42").Value = 796.1429000000001
$ws.Range("J42").Value = 69
$ws.Range("K42").Value = 2388.4287
$ws.Range("L42").Value = 207
$ws.Range("M42").Value = -2158.4287
$ws.Range("N42").Value = -667
# Row 76
$ws.Range("H76").Value = 111119464
$ws.Range("I76").Value = 200008100
$ws.Range("J76").Value = 8672
$ws.Range("K76").Value = 200008100
$ws.Range("L76").Value = 8672
$ws.Range("M76").Value = -200007785
$ws.Range("N76").Value = -9302
# Row 79
$ws.Range("H79").Value = 111119464
$ws.Range("I79").Value = 200008100
$ws.Range("J79").Value = 8672
$ws.Range("K79").Value = 200008100
$ws.Range("L79").Value = 8672
$ws.Range("M79").Value = -200007008
$ws.Range("N79").Value = -10856
# Row 92
$ws.Range("H92").Value = 805.55554
$ws.Range("I92").Value = 821.7143
$ws.Range("J92").Value = 749
$ws.Range("K92").Value = 821.7143
$ws.Range("L92").Value = 749
$ws.Range("M92").Value = 426.2857
$ws.Range("N92").Value = -3245
# Row 98
$ws.Range("H98").Value = 1363.1852
$ws.Range("I98").Value = 1372.48
$ws.Range("K98").Value = 1372.48
$ws.Range("M98").Value = 125.52
# Row 99
$ws.Range("H99").Value = 9663.333000000001
$ws.Range("I99").Value = 745
$ws.Range("K99").Value = 2235
$ws.Range("M99").Value = -737
# Row 111
$ws.Range("H111").Value = 1500
$ws.Range("J111").Value = 1500
$ws.Range("L111").Value = 4500
$ws.Range("N111").Value = -10634
# Row 113
$ws.Range("H113").Value = 40799.8
$ws.Range("I113").Value = 999.5
$ws.Range("K113").Value = 999.5
$ws.Range("M113").Value = 2254.5
# Row 122
$ws.Range("H122").Value = 1363.1852
$ws.Range("I122").Value = 1372.48
$ws.Range("K122").Value = 4117.440000000001
$ws.Range("M122").Value = -1667.440000000001
# Row 129
$ws.Range("H129").Value = 13598.75
$ws.Range("I129").Value = 758.6
$ws.Range("J129").Value = 34999
$ws.Range("K129").Value = 2275.8
$ws.Range("L129").Value = 104997
$ws.Range("M129").Value = 2724.2
$ws.Range("N129").Value = -114997
# Row 138
$ws.Range("H138").Value = 3145.75
$ws.Range("I138").Value = 937.5517
$ws.Range("J138").Value = 5211.484
$ws.Range("K138").Value = 2812.6551
$ws.Range("L138").Value = 15634.452
$ws.Range("M138").Value = 2327.3449
$ws.Range("N138").Value = -25914.452
# Row 141
$ws.Range("H141").Value = 4867.4614
$ws.Range("J141").Value = 6374.25
$ws.Range("L141").Value = 19122.75
$ws.Range("N141").Value = -29482.75

# ---- Worksheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 102
$ws.Range("H102").Value = 2786.4814
$ws.Range("I102").Value = 2864.4583
$ws.Range("J102").Value = 2162.6667
$ws.Range("K102").Value = 2864.4583
$ws.Range("L102").Value = 2162.6667
$ws.Range("M102").Value = -1242.4583
$ws.Range("N102").Value = -5406.6667
# Row 132
$ws.Range("H132").Value = 978123.0600000001
$ws.Range("I132").Value = 1686244.1
$ws.Range("J132").Value = 67681.64
$ws.Range("K132").Value = 5058732.300000001
$ws.Range("L132").Value = 203044.92
$ws.Range("M132").Value = -5056202.300000001
$ws.Range("N132").Value = -208104.92

# ---- Worksheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 9289.477999999999
$ws.Range("I99").Value = 9298.700000000001
$ws.Range("K99").Value = 9298.700000000001
$ws.Range("M99").Value = -7800.700000000001
# Row 134
$ws.Range("H134").Value = 5303444
$ws.Range("I134").Value = 7947041
$ws.Range("K134").Value = 23841123
$ws.Range("M134").Value = -23838588

# ---- Worksheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 15
$ws.Range("H15").Value = 2458.6843
$ws.Range("I15").Value = 3329.6924
$ws.Range("J15").Value = 571.5
$ws.Range("K15").Value = 3329.6924
$ws.Range("L15").Value = 571.5
$ws.Range("M15").Value = -3159.6924
$ws.Range("N15").Value = -911.5
# Row 16
$ws.Range("H16").Value = 29415412
$ws.Range("I16").Value = 71431610
$ws.Range("J16").Value = 4077.6
$ws.Range("K16").Value = 71431610
$ws.Range("L16").Value = 4077.6
$ws.Range("M16").Value = -71431323
$ws.Range("N16").Value = -4651.6
# Row 58
$ws.Range("H58").Value = 52642804
$ws.Range("I58").Value = 71436480
$ws.Range("K58").Value = 71436480
$ws.Range("M58").Value = -71436277
# Row 62
$ws.Range("H62").Value = 11648
$ws.Range("I62").Value = 13960
$ws.Range("J62").Value = 2400
$ws.Range("K62").Value = 13960
$ws.Range("L62").Value = 2400
$ws.Range("M62").Value = -13336
$ws.Range("N62").Value = -3648
# Row 65
$ws.Range("H65").Value = 11648
$ws.Range("I65").Value = 13960
$ws.Range("J65").Value = 2400
$ws.Range("K65").Value = 69800
$ws.Range("L65").Value = 12000
$ws.Range("M65").Value = -66680
$ws.Range("N65").Value = -18240
# Row 113
$ws.Range("H113").Value = 29415412
$ws.Range("I113").Value = 71431610
$ws.Range("J113").Value = 4077.6
$ws.Range("K113").Value = 71431610
$ws.Range("L113").Value = 4077.6
$ws.Range("M113").Value = -71429440
$ws.Range("N113").Value = -8417.6
# Row 132
$ws.Range("H132").Value = 7176.8096
$ws.Range("I132").Value = 6155.9473
$ws.Range("J132").Value = 16875
$ws.Range("K132").Value = 18467.8419
$ws.Range("L132").Value = 50625
$ws.Range("M132").Value = -15937.8419
$ws.Range("N132").Value = -55685
# Row 134
$ws.Range("H134").Value = 47641176
$ws.Range("I134").Value = 142872930
$ws.Range("K134").Value = 428618790
$ws.Range("M134").Value = -428616255
# Row 136
$ws.Range("H136").Value = 52642804
$ws.Range("I136").Value = 71436480
$ws.Range("K136").Value = 214309440
$ws.Range("M136").Value = -214306890

# ---- Worksheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 3508.3333
$ws.Range("I34").Value = 50
$ws.Range("J34").Value = 4200
$ws.Range("K34").Value = 150
$ws.Range("L34").Value = 12600
$ws.Range("M34").Value = -66
$ws.Range("N34").Value = -12768
# Row 52
$ws.Range("H52").Value = 3474526
$ws.Range("J52").Value = 3474526
$ws.Range("L52").Value = 10423578
$ws.Range("N52").Value = -10424110
# Row 122
$ws.Range("H122").Value = 100949.47
$ws.Range("J122").Value = 116976.23
$ws.Range("L122").Value = 1052786.07
$ws.Range("N122").Value = -1057686.07

# ---- Worksheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 4749.1787
$ws.Range("I102").Value = 3554.75
$ws.Range("J102").Value = 7735.25
$ws.Range("K102").Value = 3554.75
$ws.Range("L102").Value = 7735.25
$ws.Range("M102").Value = -1932.75
$ws.Range("N102").Value = -10979.25
# Row 132
$ws.Range("H132").Value = 25644012
$ws.Range("I132").Value = 55558220
$ws.Range("K132").Value = 166674660
$ws.Range("M132").Value = -166672130

# ---- Worksheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 4453.4
$ws.Range("I7").Value = 4059.2222
$ws.Range("J7").Value = 5044.6665
$ws.Range("K7").Value = 4059.2222
$ws.Range("L7").Value = 5044.6665
$ws.Range("M7").Value = -3947.2222
$ws.Range("N7").Value = -5268.6665
# Row 126
$ws.Range("H126").Value = 4453.4
$ws.Range("I126").Value = 4059.2222
$ws.Range("J126").Value = 5044.6665
$ws.Range("K126").Value = 12177.6666
$ws.Range("L126").Value = 15133.9995
$ws.Range("M126").Value = -9707.6666
$ws.Range("N126").Value = -20073.9995
# Row 132
$ws.Range("H132").Value = 2771.5518
$ws.Range("I132").Value = 2433.4092
$ws.Range("J132").Value = 3834.2856
$ws.Range("K132").Value = 7300.2276
$ws.Range("L132").Value = 11502.8568
$ws.Range("M132").Value = -4770.2276
$ws.Range("N132").Value = -16562.8568

# ---- Worksheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 19083.5
$ws.Range("J62").Value = 24249.666
$ws.Range("L62").Value = 24249.666
$ws.Range("N62").Value = -25497.666
# Row 65
$ws.Range("H65").Value = 19083.5
$ws.Range("J65").Value = 24249.666
$ws.Range("L65").Value = 121248.33
$ws.Range("N65").Value = -127488.33
# Row 132
$ws.Range("H132").Value = 6977.362
$ws.Range("I132").Value = 4603.2793
$ws.Range("J132").Value = 32498.75
$ws.Range("K132").Value = 13809.8379
$ws.Range("L132").Value = 97496.25
$ws.Range("M132").Value = -11279.8379
$ws.Range("N132").Value = -102556.25
